# Threat Alert Report update (2026-01-10 08:46)
#
# The two entries that used to sit in rows 4 and 5 (13-JAN-26 flynas XY-855,
# and 15-JAN-26 Nesma Airlines NE-154) are removed from the report. Deleting
# the two whole worksheet rows shifts every following row up by two, which
# reproduces the new data set (rows 6-17 become the new rows 4-15) and
# shrinks the used range from A1:K17 down to A1:K15 - exactly the same
# outcome as the published diff, without needing to retype any of the
# untouched cell values/styles by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old row 5 first, then old row 4, so the row index of the second
# deletion is unaffected by the first.
$ws.Rows.Item(5).EntireRow.Delete()
$ws.Rows.Item(4).EntireRow.Delete()
